$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF with the same style as the other header cells (A1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2 through 59: season record values for every row
for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # AD
    $ws.Cells.Item($r, 31).Value = 74   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
